$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Julio de 2020 a las 00:09"

# --- Country name re-ranking swaps (column A) caused by updated case counts ---
$ws.Range("A50").Value = "Barein"
$ws.Range("A51").Value = "Nigeria"
$ws.Range("A81").Value = "Bulgaria"
$ws.Range("A82").Value = "Estado de Palestina"
$ws.Range("A83").Value = "Noruega"
$ws.Range("A84").Value = "Senegal"
$ws.Range("A167").Value = "Comoras"
$ws.Range("A168").Value = "Isla de Man"
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4
$ws.Range("B4").Value = 4019684
$ws.Range("C4").Value = 58255
$ws.Range("D4").Value = 1884129
$ws.Range("E4").Value = 1990740
$ws.Range("G4").Value = 981
$ws.Range("H4").Value = 144815
# Row 5
$ws.Range("B5").Value = 2159654
$ws.Range("C5").Value = 38009
$ws.Range("E5").Value = 668965
$ws.Range("G5").Value = 1236
$ws.Range("H5").Value = 81487
# Row 9
$ws.Range("B9").Value = 362087
$ws.Range("C9").Value = 4406
$ws.Range("D9").Value = 248746
$ws.Range("E9").Value = 99762
$ws.Range("G9").Value = 195
$ws.Range("H9").Value = 13579
# Row 21
$ws.Range("B21").Value = 203890
$ws.Range("C21").Value = 403
$ws.Range("E21").Value = 6610
# Row 50
$ws.Range("B50").Value = 37316
$ws.Range("C50").Value = 380
$ws.Range("D50").Value = 33455
$ws.Range("E50").Value = 3732
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 129
# Row 51
$ws.Range("B51").Value = 37225
$ws.Range("D51").Value = 15333
$ws.Range("E51").Value = 21091
$ws.Range("H51").Value = 801
# Row 66
$ws.Range("B66").Value = 17881
$ws.Range("C66").Value = 567
$ws.Range("D66").Value = 9521
$ws.Range("E66").Value = 8265
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 95
# Row 81
$ws.Range("B81").Value = 9254
$ws.Range("C81").Value = 325
$ws.Range("D81").Value = 4521
$ws.Range("E81").Value = 4420
$ws.Range("G81").Value = 5
$ws.Range("H81").Value = 313
# Row 82
$ws.Range("B82").Value = 9228
$ws.Range("C82").Value = 312
$ws.Range("D82").Value = 1932
$ws.Range("E82").Value = 7232
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 64
# Row 83
$ws.Range("B83").Value = 9049
$ws.Range("C83").Value = 15
$ws.Range("D83").Value = 8138
$ws.Range("E83").Value = 656
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 255
# Row 84
$ws.Range("B84").Value = 8985
$ws.Range("C84").Value = 37
$ws.Range("D84").Value = 6044
$ws.Range("E84").Value = 2767
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 174
# Row 103
$ws.Range("B103").Value = 3817
$ws.Range("C103").Value = 69
$ws.Range("D103").Value = 2307
$ws.Range("E103").Value = 1475
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 35
# Row 125
$ws.Range("B125").Value = 1954
$ws.Range("C125").Value = 5
$ws.Range("E125").Value = 1125
# Row 131
$ws.Range("B131").Value = 1655
$ws.Range("C131").Value = 26
$ws.Range("D131").Value = 848
$ws.Range("E131").Value = 802
# Row 138
$ws.Range("D138").Value = 1034
$ws.Range("E138").Value = 178
# Row 151
$ws.Range("B151").Value = 790
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 560
$ws.Range("E151").Value = 215
# Row 167
$ws.Range("B167").Value = 337
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 319
$ws.Range("E167").Value = 11
$ws.Range("H167").Value = 7
# Row 168
$ws.Range("B168").Value = 336
$ws.Range("D168").Value = 312
$ws.Range("E168").Value = 0
$ws.Range("H168").Value = 24
